# Update works images 2026-01-21 15:35:00
#
# "Sheet1" holds the pool of still-unused generated names (one per row,
# column A). "used" is the append-only log of names that have been
# consumed, together with the source image filename and the timestamp
# it was used.
#
# This run consumes the next available name ("3cs680lw", row 2 of
# Sheet1) for a newly generated image, so:
#   1. remove it from the front of the Sheet1 pool (row 2 is deleted,
#      shifting every following row up by one), and
#   2. append it as a new row at the bottom of the "used" log with the
#      image's filename and usage timestamp.

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# 1) Pop the consumed name off the top of the pool (row 1 is the header-less
#    first entry; the name being consumed lives in row 2).
$namesSheet.Rows(2).Delete()

# 2) Log it as used, in the first empty row after the existing log.
$nextRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row + 1

$usedSheet.Cells.Item($nextRow, 1).Value = "3cs680lw"
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月21日 15_24_34.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-21 15:34:54"
